# Preparing final SEMAINE submission.
# Add new SVM-geometry result rows to the SEMAINE sheet, re-order a couple of
# pre-existing rows, and leave the UI focused on the SEMAINE tab / the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SEMAINE")

# --- Row 29: new "DISFA SVM geometry stat" result row -----------------------
$ws.Range("A29").Value = "DISFA SVM geometry stat"
$ws.Range("B29").Value = 0.21
$ws.Range("C29").Value = 0.64
$ws.Range("D29").Value = 0.32
$ws.Range("E29").Value = 0.38
$ws.Range("F29").Value = 0.76
$ws.Range("G29").Value = 0.51
$ws.Range("H29").Value = 0.1
$ws.Range("I29").Value = 0.88
$ws.Range("J29").Value = 0.18
$ws.Range("K29").Value = 0.3
$ws.Range("L29").Value = 0.76
$ws.Range("M29").Value = 0.43

# --- Row 30: "DISFA SVM geometry dyn" (moved down from old row 29, new data) -
$ws.Range("A30").Value = "DISFA SVM geometry dyn"
$ws.Range("B30").Value = 0.603
$ws.Range("C30").Value = 0.474
$ws.Range("D30").Value = 0.531
$ws.Range("E30").Value = 0.741
$ws.Range("F30").Value = 0.176
$ws.Range("G30").Value = 0.284
$ws.Range("H30").Value = 0.158
$ws.Range("I30").Value = 0.715
$ws.Range("J30").Value = 0.259
$ws.Range("K30").Value = 0.278
$ws.Range("L30").Value = 0.405
$ws.Range("M30").Value = 0.33

# --- Row 31: "BP4D SVM geometry dyn" (moved down from old row 30) -----------
$ws.Range("A31").Value = "BP4D SVM geometry dyn"
$ws.Range("B31").Value = 0.2652
$ws.Range("C31").Value = 0.547
$ws.Range("D31").Value = 0.3572
$ws.Range("E31").Value = 0.3364
$ws.Range("F31").Value = 0.8099
$ws.Range("G31").Value = 0.4754
$ws.Range("H31").Value = 0.0609
$ws.Range("I31").Value = 0.9485
$ws.Range("J31").Value = 0.1145

# --- Row 32: "BP4D SVM geometry stat" (moved down from old row 31) ----------
$ws.Range("A32").Value = "BP4D SVM geometry stat"
$ws.Range("B32").Value = 0.1608
$ws.Range("C32").Value = 0.8957
$ws.Range("D32").Value = 0.2727
$ws.Range("E32").Value = 0.3331
$ws.Range("F32").Value = 0.9528
$ws.Range("G32").Value = 0.4936
$ws.Range("H32").Value = 0.0655
$ws.Range("I32").Value = 0.921
$ws.Range("J32").Value = 0.1223

# --- Row 33: new "Combined SVM geometry stat" result row --------------------
$ws.Range("A33").Value = "Combined SVM geometry stat"
$ws.Range("B33").Value = 0.3035
$ws.Range("C33").Value = 0.588
$ws.Range("D33").Value = 0.4004

# --- Row 34: new "Combined SVM geometry dyn v sem" result row --------------
$ws.Range("A34").Value = "Combined SVM geometry dyn v sem"
$ws.Range("B34").Value = 0.63
$ws.Range("C34").Value = 0.5153
$ws.Range("D34").Value = 0.5669
$ws.Range("E34").Value = 0.4251
$ws.Range("F34").Value = 0.6187
$ws.Range("G34").Value = 0.5039
$ws.Range("H34").Value = 0.1922
$ws.Range("I34").Value = 0.2468
$ws.Range("J34").Value = 0.2161
$ws.Range("K34").Value = 0.2881
$ws.Range("L34").Value = 0.5751
$ws.Range("M34").Value = 0.3838

# --- Row 35: fresh blank templated row (copy the B:V formatting down one row)
$ws.Range("B34:V34").Copy()
$ws.Range("B35:V35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Leave the UI pointed at the SEMAINE tab, with the new row selected ----
$wsBP4D = $wb.Worksheets.Item("BP4D")
$wsBP4D.Activate()
$wsBP4D.Range("J27").Select()

$ws.Activate()
$ws.Range("D34").Select()
